# Actualización desde MV -datos-
# Append 5 new daily rows (02-10-2021 .. 06-10-2021) to the bottom of the
# "Reinversión de bonos bancarios" daily table, following the same pattern
# as the existing rows (B=449, C=0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 249
$dates = @("02-10-2021", "03-10-2021", "04-10-2021", "05-10-2021", "06-10-2021")

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Use a formula that evaluates to the literal date string so Excel
    # stores it as text (a "dd-mm-yyyy" string would otherwise be
    # auto-parsed as a date), then flatten it to a static value via
    # copy/paste-special so no formula or extra cell style remains.
    $ws.Cells.Item($row, 1).Formula = '="' + $dates[$i] + '"'

    $ws.Cells.Item($row, 2).Value = 449
    $ws.Cells.Item($row, 3).Value = 0
}

$dataRange = $ws.Range("A" + $startRow + ":A" + ($startRow + $dates.Length - 1))
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)
